# Commit: "Changed GPA as number"
#
# Column X ("GPA") held the literal text "--" for students whose GPA was
# not applicable/available. Re-enter those cells as the number 0 instead
# of the text placeholder so the column is numeric throughout.
#
# Rows below are exactly the ones whose column-X value is currently the
# text "--".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    2,3,4,7,8,9,11,12,13,14,15,17,20,22,23,24,25,31,32,33,34,37,38,39,40,
    43,44,45,47,48,51,53,55,58,59,62,63,66,70,71,73,74,76,78,79,80,81,82,
    83,84,85,86,87,88,90,92,94,95,96,98,99,100,101,102,104,106,107,108,
    110,111,113,114,118,120,121,124,125,126,129,130,131,132,133,135,136,
    137,138,139,140,141,142,143,145,148,149,150,152,153,154,155,156,157,
    161,162,163,164,165,166,167,168,169,173
)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 24).Value = 0
}

# Keep the active selection consistent with the edited column.
$ws.Range("X2").Select()
